$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02264433333333333
$ws.Range("N2").Value = 0.06793299999999999
$ws.Range("O2").Value = 0.08454793804489194
$ws.Range("P2").Value = 0.08454793804489193
$ws.Range("Q2").Value = 3.278521419521666
$ws.Range("R2").Value = 29.50669277569499
$ws.Range("S2").Value = 0.02054554069939457
$ws.Range("T2").Value = 0.02122260470439368

# Row 3
$ws.Range("G3").Value = 144.783305
$ws.Range("H3").Value = 434.349915
$ws.Range("I3").Value = 0.2430046335191003
$ws.Range("J3").Value = 0.251012682214973
$ws.Range("M3").Value = 0.245184
$ws.Range("N3").Value = 0.735552
$ws.Range("O3").Value = 0.915452061955108
$ws.Range("P3").Value = 0.9154520619551081
$ws.Range("Q3").Value = 35.49854985312
$ws.Range("R3").Value = 319.48694867808
$ws.Range("S3").Value = 0.2224590928197057
$ws.Range("T3").Value = 0.2297900775105793

# Row 4
$ws.Range("G4").Value = 82.24887099999999
$ws.Range("I4").Value = 0.1380466950572427
$ws.Range("J4").Value = 0.1425959278859072
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02264433333333333
$ws.Range("N4").Value = 0.06793299999999999
$ws.Range("O4").Value = 0.08454793804489194
$ws.Range("P4").Value = 0.08454793804489193
$ws.Range("Q4").Value = 1.862470851214333
$ws.Range("R4").Value = 16.762237660929
$ws.Range("S4").Value = 0.01167156342100184
$ws.Range("T4").Value = 0.01205619167635156

# Row 5
$ws.Range("G5").Value = 82.24887099999999
$ws.Range("I5").Value = 0.1380466950572427
$ws.Range("J5").Value = 0.1425959278859072
$ws.Range("M5").Value = 0.245184
$ws.Range("N5").Value = 0.735552
$ws.Range("O5").Value = 0.915452061955108
$ws.Range("P5").Value = 0.9154520619551081
$ws.Range("Q5").Value = 20.166107187264
$ws.Range("R5").Value = 181.494964685376
$ws.Range("S5").Value = 0.1263751316362408
$ws.Range("T5").Value = 0.1305397362095556

# Row 6
$ws.Range("G6").Value = 163.8590903333333
$ws.Range("H6").Value = 491.577271
$ws.Range("I6").Value = 0.2750214756820535
$ws.Range("J6").Value = 0.284084617144743
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02264433333333333
$ws.Range("N6").Value = 0.06793299999999999
$ws.Range("O6").Value = 0.08454793804489194
$ws.Range("P6").Value = 0.08454793804489193
$ws.Range("Q6").Value = 3.710479861204778
$ws.Range("R6").Value = 33.394318750843
$ws.Range("S6").Value = 0.02325249868698101
$ws.Range("T6").Value = 0.02401876860986057

# Row 7
$ws.Range("G7").Value = 163.8590903333333
$ws.Range("H7").Value = 491.577271
$ws.Range("I7").Value = 0.2750214756820535
$ws.Range("J7").Value = 0.284084617144743
$ws.Range("M7").Value = 0.245184
$ws.Range("N7").Value = 0.735552
$ws.Range("O7").Value = 0.915452061955108
$ws.Range("P7").Value = 0.9154520619551081
$ws.Range("Q7").Value = 40.175627204288
$ws.Range("R7").Value = 361.580644838592
$ws.Range("S7").Value = 0.2517689769950724
$ws.Range("T7").Value = 0.2600658485348824

# Row 8
$ws.Range("G8").Value = 57.0238095
$ws.Range("H8").Value = 114.047619
$ws.Range("I8").Value = 0.09570889357312636
$ws.Range("J8").Value = 0.06590860906562239
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02264433333333333
$ws.Range("N8").Value = 0.06793299999999999
$ws.Range("O8").Value = 0.08454793804489194
$ws.Range("P8").Value = 0.08454793804489193
$ws.Range("Q8").Value = 1.2912661502545
$ws.Range("R8").Value = 7.747596901526999
$ws.Range("S8").Value = 0.008091989604165843
$ws.Range("T8").Value = 0.005572436995905244

# Row 9
$ws.Range("G9").Value = 57.0238095
$ws.Range("H9").Value = 114.047619
$ws.Range("I9").Value = 0.09570889357312636
$ws.Range("J9").Value = 0.06590860906562239
$ws.Range("M9").Value = 0.245184
$ws.Range("N9").Value = 0.735552
$ws.Range("O9").Value = 0.915452061955108
$ws.Range("P9").Value = 0.9154520619551081
$ws.Range("Q9").Value = 13.981325708448
$ws.Range("R9").Value = 83.887954250688
$ws.Range("S9").Value = 0.08761690396896051
$ws.Range("T9").Value = 0.06033617206971715

# Row 10
$ws.Range("G10").Value = 147.8896333333333
$ws.Range("H10").Value = 443.6689
$ws.Range("I10").Value = 0.2482183021684772
$ws.Range("J10").Value = 0.2563981636887546
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02264433333333333
$ws.Range("N10").Value = 0.06793299999999999
$ws.Range("O10").Value = 0.08454793804489194
$ws.Range("P10").Value = 0.08454793804489193
$ws.Range("Q10").Value = 3.348862153744444
$ws.Range("R10").Value = 30.1397593837
$ws.Range("S10").Value = 0.02098634563334868
$ws.Range("T10").Value = 0.02167793605838088

# Row 11
$ws.Range("G11").Value = 147.8896333333333
$ws.Range("H11").Value = 443.6689
$ws.Range("I11").Value = 0.2482183021684772
$ws.Range("J11").Value = 0.2563981636887546
$ws.Range("M11").Value = 0.245184
$ws.Range("N11").Value = 0.735552
$ws.Range("O11").Value = 0.915452061955108
$ws.Range("P11").Value = 0.9154520619551081
$ws.Range("Q11").Value = 36.2601718592
$ws.Range("R11").Value = 326.3415467328
$ws.Range("S11").Value = 0.2272319565351285
$ws.Range("T11").Value = 0.2347202276303737
